# Apply the LinuxForHealth re-brand + StructureDefinition metadata update.
$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL (B2): ibm.com -> linuxforhealth.org
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/child-organization-hierarchy-level-description"

# Version (B3): 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date (B8): 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher (B9): Alvearie Team -> LinuxForHealth Team
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet ---------------------------------------------------
# The "Constraint(s)" text (ele-1 / ext-1) was mis-placed on the
# `Extension` root row (row 2) instead of the `Extension.extension`
# row (row 4). Move it down to row 4 and clear the stray copy on row 2.
$elements = $wb.Worksheets.Item("Elements")

$constraintText = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

$elements.Range("AI2").Value = ""
$elements.Range("AI4").Value = $constraintText

# Fixed Value (Q5) on the Extension.url row repeats the same URL - keep in sync.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/child-organization-hierarchy-level-description"
